# Task now has a duration instead of a start time: the "minutes" column
# header becomes "duration in minutes", and the help text in the
# instructions text box is updated to match (task_name must be unique /
# minutes = duration of the task instead of a start offset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header change: B1 "minutes" -> "duration in minutes" -----------------
$ws.Range("B1").Value = "duration in minutes"

# Put the active selection on B1 (matches the saved cursor position).
[void]$ws.Range("B1").Select()

# --- Update the "HowTo" help text box on the sheet ------------------------
if ($ws.Shapes.Count -ge 1) {
    $shp = $ws.Shapes.Item(1)

    $helpText = @"
HowTo
task_name: The name of the task. Must be unique

minutes. How long the task takes.

command: Command to execute (Leav empty if you do not want to execute a command)

title: Title of the text screen (leave empty if you are using a command)
Description (Optional to title): Description

Variables:
{id} -> Participant ID
{experiment} -> Experiment Name
{startTime} -> Time when the experiment started
{timestamp} -> Current time stamp. Format: YYYY.mm.dd hh:mm:ss
{scriptCount} -> A counter that increments with each execution of the command within a task (used conly for commands)
"@

    $shp.TextFrame.Characters().Text = $helpText
}
